$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 10
$ws.Range("A10").Value = 111577193
$ws.Range("B10").Value = 89845
$ws.Range("E10").Value = 1209
$ws.Range('F10').Value = 'Rynkskinn'
$ws.Range('G10').Value = 'Phlebia centrifuga'
$ws.Range('H10').Value = 'P.Karst.'
$ws.Range("Q10").Value = 562885.3077477051
$ws.Range("R10").Value = 6954427.514711756
$ws.Range('Z10').Value = '17:24'
$ws.Range('AB10').Value = '17:24'

# Row 11
$ws.Range("A11").Value = 111577080
$ws.Range("B11").Value = 89405
$ws.Range('D11').Value = 'NT'
$ws.Range("E11").Value = 1202
$ws.Range('F11').Value = 'Ullticka'
$ws.Range('G11').Value = 'Phellinidium ferrugineofuscum'
$ws.Range('H11').Value = '(P.Karst.) Fiasson & Niemelä'
$ws.Range("Q11").Value = 562937.8525077751
$ws.Range("R11").Value = 6954467.524316943
$ws.Range('Z11').Value = '17:24'
$ws.Range('AB11').Value = '17:24'

# Row 12
$ws.Range("A12").Value = 111577838
$ws.Range("B12").Value = 96348
$ws.Range("E12").Value = 220787
$ws.Range('F12').Value = 'Knärot'
$ws.Range('G12').Value = 'Goodyera repens'
$ws.Range('H12').Value = '(L.) R. Br.'
$ws.Range("Q12").Value = 562801.9298062191
$ws.Range("R12").Value = 6954389.67147268
$ws.Range('Z12').Value = '18:14'
$ws.Range('AB12').Value = '18:14'

# Row 13
$ws.Range("A13").Value = 111577611
$ws.Range("Q13").Value = 562810.9079359611
$ws.Range("R13").Value = 6954400.856378952
$ws.Range('Z13').Value = '00:00'
$ws.Range('AB13').Value = '00:00'

# Row 15
$ws.Range("A15").Value = 111577804
$ws.Range("B15").Value = 96348
$ws.Range('D15').Value = 'VU'
$ws.Range("E15").Value = 220787
$ws.Range('F15').Value = 'Knärot'
$ws.Range('G15').Value = 'Goodyera repens'
$ws.Range('H15').Value = '(L.) R. Br.'
$ws.Range("Q15").Value = 562814.775380839
$ws.Range("R15").Value = 6954390.834027934
$ws.Range('Z15').Value = '00:00'
$ws.Range('AB15').Value = '00:00'

# Row 16
$ws.Range("A16").Value = 111577994
$ws.Range("B16").Value = 89686
$ws.Range('D16').Value = 'NT'
$ws.Range("E16").Value = 658
$ws.Range('F16').Value = 'Rosenticka'
$ws.Range('G16').Value = 'Rhodofomes roseus'
$ws.Range('H16').Value = '(Alb. & Schwein.) Kotl. & Pouzar'
$ws.Range("Q16").Value = 562890.6102569005
$ws.Range("R16").Value = 6954486.814324431

# Row 18
$ws.Range("A18").Value = 111575408
$ws.Range("Q18").Value = 562539.034657649
$ws.Range("R18").Value = 6954609.073577877
$ws.Range('Z18').Value = '16:21'
$ws.Range('AB18').Value = '16:21'

# Row 19
$ws.Range("A19").Value = 111578090
$ws.Range("B19").Value = 89405
$ws.Range("E19").Value = 1202
$ws.Range('F19').Value = 'Ullticka'
$ws.Range('G19').Value = 'Phellinidium ferrugineofuscum'
$ws.Range('H19').Value = '(P.Karst.) Fiasson & Niemelä'
$ws.Range("Q19").Value = 562890.6102569005
$ws.Range("R19").Value = 6954486.814324431
$ws.Range('Z19').Value = '18:30'
$ws.Range('AB19').Value = '18:30'

# Row 20
$ws.Range("A20").Value = 111577743
$ws.Range("Q20").Value = 562802.8660743404
$ws.Range("R20").Value = 6954388.771485241
$ws.Range('Z20').Value = '17:59'
$ws.Range('AB20').Value = '17:59'

# Row 21
$ws.Range("A21").Value = 111577591
$ws.Range("B21").Value = 56543
$ws.Range("E21").Value = 103021
$ws.Range('F21').Value = 'Talltita'
$ws.Range('G21').Value = 'Poecile montanus'
$ws.Range('H21').Value = '(Conrad von Baldenstein, 1827)'
$ws.Range('M21').Value = 'lockläte, övriga läten'
$ws.Range("Q21").Value = 562822.1033927511
$ws.Range("R21").Value = 6954368.028004575
$ws.Range('Z21').Value = '00:00'
$ws.Range('AB21').Value = '00:00'

# Row 22
$ws.Range("A22").Value = 111577347
$ws.Range("B22").Value = 96348
$ws.Range('D22').Value = 'VU'
$ws.Range("E22").Value = 220787
$ws.Range('F22').Value = 'Knärot'
$ws.Range('G22').Value = 'Goodyera repens'
$ws.Range('H22').Value = '(L.) R. Br.'
$ws.Range("M22").ClearContents() | Out-Null
$ws.Range("Q22").Value = 562796.503171768
$ws.Range("R22").Value = 6954336.792844097

# Row 23
$ws.Range("A23").Value = 111577180
$ws.Range("B23").Value = 89686
$ws.Range("E23").Value = 658
$ws.Range('F23').Value = 'Rosenticka'
$ws.Range('G23').Value = 'Rhodofomes roseus'
$ws.Range('H23').Value = '(Alb. & Schwein.) Kotl. & Pouzar'
$ws.Range("Q23").Value = 562883.198926247
$ws.Range("R23").Value = 6954441.700568204
$ws.Range('Z23').Value = '17:24'
$ws.Range('AB23').Value = '17:24'
